# Applies the "Fitting probability distribution to results and plotting" edit
# to the Inputs worksheet of the PSI_inputs workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs")

# No_cycles row (row 12): was a single "BE" value of 10 in column B.
# Now split across LE/BE/HE columns: 1 / 10 / 100
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 100

# Output_dist (row 15): distribution changed from Normal to Weibull
$ws.Range("B15").Value = "Weibull"

# Young's modulus "E" row (25) replaced with bending stiffness "EI" row
$ws.Range("A25").Value = "EI"
$ws.Range("B25").Value = 61675
$ws.Range("C25").Value = 61675
$ws.Range("D25").Value = 61675
$ws.Range("G25").Value = "% pipe bending stiffness (kNm2)"

# Restore view state (selection / scroll position) on the Inputs sheet
$ws.Activate() | Out-Null
$ws.Range("B16").Select() | Out-Null
